# Update the crop list in column A to reflect corrected/added crop codes.
# - Added "corn- ornamental, glass gem" (alphabetically after "cilantro")
# - Added "pepper- califonia wonder" (alphabetically after "pepper- anaheim")
# - Renamed "pepper- jalapeno- early" -> "pepper- jalapeno, early"
# - Removed "pepper- sweet- ca wonder"
# - Renamed "pumpkin- ct field" -> "pumpkin- connecticut field"
# - Removed "sweet corn- ornamental- glas gm"
# The resulting list keeps the same row count (23 crops), each row's
# crop_code (column B) stays the same 1..23 sequence already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$crops = @(
    "arrugula- garden tangy",
    "basil- sweet",
    "bean- fin de bagnol",
    "carrot- touchon",
    "cilantro",
    "corn- ornamental, glass gem",
    "dill- hera",
    "eggplant- black beauty",
    "garlic- german extra hardy hardneck",
    "kale- lacinato",
    "lettuce- gourmet blend",
    "pepper- anaheim",
    "pepper- califonia wonder",
    "pepper- ghost ",
    "pepper- habanero",
    "pepper- jalapeno, early",
    "pepper- sweet banana",
    "pumpkin- connecticut field",
    "sage- broad leaf",
    "sunflower- mammoth",
    "sweet corn- country gentleman",
    "tomato- amish",
    "tomato- roma"
)

for ($i = 0; $i -lt $crops.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $crops[$i]
}
